$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (random_forest) - updated metrics
$ws.Range("B2").Value = 4.4462543499490916
$ws.Range("C2").Value = 0.37811088586590397
$ws.Range("D2").Value = 3.5804165403231774
$ws.Range("E2").Value = 0.35973861114165051
$ws.Range("F2").Value = 0.59978213639758438
$ws.Range("G2").Value = 0.54834769221899571
$ws.Range("H2").Value = 0.64026138885834949
$ws.Range("I2").Value = 0.8325092811520225

# Row 3 (lsboost) - updated metrics
$ws.Range("B3").Value = 4.6637581786857902
$ws.Range("C3").Value = 0.39660748072755436
$ws.Range("D3").Value = 3.654055095483784
$ws.Range("E3").Value = 0.39579517427749217
$ws.Range("F3").Value = 0.62912254313249039
$ws.Range("G3").Value = 0.55962557883523145
$ws.Range("H3").Value = 0.60420482572250789
$ws.Range("I3").Value = 0.78756006301541914

# Row 4 (neural_network) - updated metrics
$ws.Range("B4").Value = 4.1921153396111563
$ws.Range("C4").Value = 0.35649882349410983
$ws.Range("D4").Value = 3.380758035404825
$ws.Range("E4").Value = 0.31979001076903546
$ws.Range("F4").Value = 0.5654997884783296
$ws.Range("G4").Value = 0.51776960747079182
$ws.Range("H4").Value = 0.68020998923096454
$ws.Range("I4").Value = 0.85226031249339207

# Row 5 (old_model) values remain unchanged
